# Fruta / hortaliza, semanal
# A new weekly record is inserted as the first data row of the "Femacal de
# La Calera - Papaya" block (row 111), pushing the previously existing rows
# 111-124 down to 112-125.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 111; this shifts rows 111-124 down to 112-125 and
# carries over formatting (e.g. the date style on column D) automatically.
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row 111 with the new weekly record.
$ws.Range("A111").Value = 3
$ws.Range("B111").Value = "Femacal de La Calera"
$ws.Range("C111").Value = "Coquimbo"
$ws.Range("D111").Value = "2023-10-13"
$ws.Range("E111").Value = 5
$ws.Range("F111").Value = "Fruta"
$ws.Range("G111").Value = 100108
$ws.Range("H111").Value = "Tropicales y subtropicales"
$ws.Range("I111").Value = 100108004
$ws.Range("J111").Value = "Papaya"
$ws.Range("K111").Value = "Cultivar IV Región"
$ws.Range("L111").Value = "Primera"
$ws.Range("M111").Value = 96
$ws.Range("N111").Value = 16000
$ws.Range("O111").Value = 17000
$ws.Range("P111").Value = 16417
$ws.Range("Q111").Value = "$/bandeja 10 kilos"
$ws.Range("R111").Value = "Provincia del Elquí"
$ws.Range("S111").Value = 1642
$ws.Range("T111").Value = 10
